$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per latest scrape
# Force text number format so numeric-looking strings (e.g. "592.14", "1.00")
# are preserved exactly as text rather than being reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.740.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.64%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.520.77"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.79"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.976.68"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.552.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.513.90"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.96"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.644.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0990"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.75"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.66"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.23%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.562"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "146.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0279"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.64%  "
